$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Update a few source values on Sheet1 (dependent formulas recalc automatically) ---
$ws1.Range("J6").Value = 56
$ws1.Range("K6").Value = 100
$ws1.Range("L6").Value = 16

# --- Add new worksheet "Sheet2" positioned right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- Populate Sheet2 data ---
# Columns A-C block (quantity / unit price / subtotal)
$ws2.Range("A1").Value = 5
$ws2.Range("B1").Value = 250
$ws2.Range("C1").Formula = "=B1*A1"

$ws2.Range("A2").Value = 15
$ws2.Range("B2").Value = 16
$ws2.Range("C2").Formula = "=B2*A2"

$ws2.Range("A3").Value = 30
$ws2.Range("B3").Value = 21
$ws2.Range("C3").Formula = "=B3*A3"

$ws2.Range("A4").Value = 60
$ws2.Range("B4").Value = 39
$ws2.Range("C4").Formula = "=B4*A4"

$ws2.Range("A5").Value = 120
$ws2.Range("B5").Value = 14
$ws2.Range("C5").Formula = "=B5*A5"

$ws2.Range("A6").Value = 480
$ws2.Range("B6").Value = 14
$ws2.Range("C6").Formula = "=B6*A6"

# Columns D-F block (quantity / unit price / subtotal)
$ws2.Range("D1").Value = 5
$ws2.Range("E1").Value = 74
$ws2.Range("F1").Formula = "=E1*D1"

$ws2.Range("D2").Value = 15
$ws2.Range("E2").Value = 1
$ws2.Range("F2").Formula = "=E2*D2"

$ws2.Range("D3").Value = 30
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Formula = "=E3*D3"

$ws2.Range("D4").Value = 60
$ws2.Range("E4").Value = 43
$ws2.Range("F4").Formula = "=E4*D4"

$ws2.Range("D5").Value = 120
$ws2.Range("E5").Value = 18
$ws2.Range("F5").Formula = "=E5*D5"

$ws2.Range("D6").Value = 480
$ws2.Range("E6").Value = 3
$ws2.Range("F6").Formula = "=E6*D6"

# Totals row
$ws2.Range("C7").Formula = "=SUM(C1:C6)"
$ws2.Range("F7").Formula = "=SUM(F1:F6)"
$ws2.Range("G7").Formula = "=C7+F7"

# Derived rows
$ws2.Range("G8").Formula = "=G7/60"
$ws2.Range("G9").Formula = "=G8/24"

# --- Sheet2 view state (selection cursor left at G10) ---
$ws2.Range("G10").Select() | Out-Null

# --- Re-activate Sheet1 so it remains the selected/visible tab ---
$ws1.Activate()
$ws1.Range("M12").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save() | Out-Null
